$wb = $excel.ActiveWorkbook

# Updated correlation matrix for the "Imports" sheet (B2:G7)
$importsValues = @(
    @(1, 0.493367024071237, 0.659662165314355, 0.770933069422289, 0.549502247347733, 0.217973605171083),
    @(0.493367024071237, 1, 0.551897196626651, 0.447316253062892, -0.0655443226404276, -0.298608174491895),
    @(0.659662165314355, 0.551897196626651, 1, 0.78361925100665, 0.392488167211744, -0.12359513834032),
    @(0.770933069422289, 0.447316253062892, 0.78361925100665, 1, 0.55063452664881, -0.22625033920859),
    @(0.549502247347733, -0.0655443226404276, 0.392488167211744, 0.55063452664881, 1, 0.0668730984730335),
    @(0.217973605171083, -0.298608174491895, -0.12359513834032, -0.22625033920859, 0.0668730984730335, 1)
)

# Updated correlation matrix for the "Exports" sheet (B2:G7)
$exportsValues = @(
    @(1, 0.497973550996671, -0.0643901413606171, 0.463761020727406, -0.331023633136514, 0.542387538504438),
    @(0.497973550996671, 1, 0.190702768235947, 0.508427566507527, 0.505597673361365, 0.0733371328413457),
    @(-0.0643901413606171, 0.190702768235947, 1, 0.0458749036578238, 0.0737230234435839, -0.0448196890732375),
    @(0.463761020727406, 0.508427566507527, 0.0458749036578238, 1, 0.262866408852198, 0.0447915176767876),
    @(-0.331023633136514, 0.505597673361365, 0.0737230234435839, 0.262866408852198, 1, -0.434127752053819),
    @(0.542387538504438, 0.0733371328413457, -0.0448196890732375, 0.0447915176767876, -0.434127752053819, 1)
)

$wsImports = $wb.Worksheets.Item("Imports")
$wsExports = $wb.Worksheets.Item("Exports")

for ($r = 0; $r -lt 6; $r++) {
    for ($c = 0; $c -lt 6; $c++) {
        $wsImports.Cells.Item($r + 2, $c + 2).Value = $importsValues[$r][$c]
        $wsExports.Cells.Item($r + 2, $c + 2).Value = $exportsValues[$r][$c]
    }
}
